$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Harshal Patel"

# Make sure the whole target range is treated as text, matching the
# source data (all values are stored as strings, even numeric-looking ones)
$ws.Range("A1:M9").NumberFormat = "@"

# Header row (a new "matchNo" column is inserted before the existing columns,
# shifting every other header one column to the right)
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# Row 2
$ws.Range("A2").Value = "Eliminator"
$ws.Range("B2").Value = "Royal Challengers Bangalore"
$ws.Range("C2").Value = "Harshal Patel"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "8"
$ws.Range("F2").Value = "6"
$ws.Range("G2").Value = "1"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "133.33"
$ws.Range("J2").Value = "Kolkata Knight Riders"
$ws.Range("K2").Value = "Sharjah"
$ws.Range("L2").Value = "October 11"
$ws.Range("M2").Value = "KKR won by 4 wickets (with 2 balls remaining)"

# Row 3
$ws.Range("A3").Value = "31st"
$ws.Range("B3").Value = "Royal Challengers Bangalore"
$ws.Range("C3").Value = "Harshal Patel"
$ws.Range("D3").Value = "b Ferguson"
$ws.Range("E3").Value = "12"
$ws.Range("F3").Value = "10"
$ws.Range("G3").Value = "2"
$ws.Range("H3").Value = "0"
$ws.Range("I3").Value = "120.00"
$ws.Range("J3").Value = "Kolkata Knight Riders"
$ws.Range("K3").Value = "Abu Dhabi"
$ws.Range("L3").Value = "September 20"
$ws.Range("M3").Value = "KKR won by 9 wickets (with 60 balls remaining)"

# Row 4
$ws.Range("A4").Value = "19th"
$ws.Range("B4").Value = "Royal Challengers Bangalore"
$ws.Range("C4").Value = "Harshal Patel"
$ws.Range("D4").Value = "b Imran Tahir"
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "8"
$ws.Range("G4").Value = "0"
$ws.Range("H4").Value = "0"
$ws.Range("I4").Value = "0.00"
$ws.Range("J4").Value = "Chennai Super Kings"
$ws.Range("K4").Value = "Wankhede"
$ws.Range("L4").Value = "April 25"
$ws.Range("M4").Value = "Super Kings won by 69 runs"

# Row 5
$ws.Range("A5").Value = "48th"
$ws.Range("B5").Value = "Royal Challengers Bangalore"
$ws.Range("C5").Value = "Harshal Patel"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "1"
$ws.Range("G5").Value = "0"
$ws.Range("H5").Value = "0"
$ws.Range("I5").Value = "100.00"
$ws.Range("J5").Value = "Punjab Kings"
$ws.Range("K5").Value = "Sharjah"
$ws.Range("L5").Value = "October 03"
$ws.Range("M5").Value = "RCB won by 6 runs"

# Row 6
$ws.Range("A6").Value = "35th"
$ws.Range("B6").Value = "Royal Challengers Bangalore"
$ws.Range("C6").Value = "Harshal Patel"
$ws.Range("D6").Value = "c Raina b Bravo"
$ws.Range("E6").Value = "3"
$ws.Range("F6").Value = "5"
$ws.Range("G6").Value = "0"
$ws.Range("H6").Value = "0"
$ws.Range("I6").Value = "60.00"
$ws.Range("J6").Value = "Chennai Super Kings"
$ws.Range("K6").Value = "Sharjah"
$ws.Range("L6").Value = "September 24"
$ws.Range("M6").Value = "Super Kings won by 6 wickets (with 11 balls remaining)"

# Row 7
$ws.Range("A7").Value = "26th"
$ws.Range("B7").Value = "Royal Challengers Bangalore"
$ws.Range("C7").Value = "Harshal Patel"
$ws.Range("D7").Value = "c Ravi Bishnoi b Mohammed Shami"
$ws.Range("E7").Value = "31"
$ws.Range("F7").Value = "13"
$ws.Range("G7").Value = "3"
$ws.Range("H7").Value = "2"
$ws.Range("I7").Value = "238.46"
$ws.Range("J7").Value = "Punjab Kings"
$ws.Range("K7").Value = "Ahmedabad"
$ws.Range("L7").Value = "April 30"
$ws.Range("M7").Value = "Punjab Kings won by 34 runs"

# Row 8
$ws.Range("A8").Value = "6th"
$ws.Range("B8").Value = "Royal Challengers Bangalore"
$ws.Range("C8").Value = "Harshal Patel"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = "0"
$ws.Range("F8").Value = "0"
$ws.Range("G8").Value = "0"
$ws.Range("H8").Value = "0"
$ws.Range("I8").Value = "-"
$ws.Range("J8").Value = "Sunrisers Hyderabad"
$ws.Range("K8").Value = "Chennai"
$ws.Range("L8").Value = "April 14"
$ws.Range("M8").Value = "RCB won by 6 runs"

# Row 9
$ws.Range("A9").Value = "1st"
$ws.Range("B9").Value = "Royal Challengers Bangalore"
$ws.Range("C9").Value = "Harshal Patel"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = "4"
$ws.Range("F9").Value = "3"
$ws.Range("G9").Value = "0"
$ws.Range("H9").Value = "0"
$ws.Range("I9").Value = "133.33"
$ws.Range("J9").Value = "Mumbai Indians"
$ws.Range("K9").Value = "Chennai"
$ws.Range("L9").Value = "April 09"
$ws.Range("M9").Value = "RCB won by 2 wickets"
